$d = $word.ActiveDocument

function Get-XmlPkg([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Remove the stray empty run (<w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>)
#    that trails the "Home Station" run in paragraph 1. The empty run has no
#    addressable characters, so we must replace a range that spans into the
#    following paragraph (paragraph 2) to force the engine to drop it, then
#    re-emit paragraph 2's content unchanged.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p3 = $d.Paragraphs(3)
$rng1 = $d.Range($p1.Range.Start, $p3.Range.Start)
$body1 = '<w:body>' + `
    '<w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b w:val="1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:b w:val="1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Home Station</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Petra (neutral smiling): It’s a lot nicer today than it was yesterday, huh?</w:t></w:r></w:p>' + `
    '</w:body>'
$rng1.InsertXML((Get-XmlPkg $body1))

# ---------------------------------------------------------------------------
# 2) Change "raised_eyebrow" -> "skeptical" in the "Habit?" line, and insert a
#    new paragraph "Petra (neutral sigh):" right after it.
# ---------------------------------------------------------------------------
$habitParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Habit? Who have you been walking home*") {
        $habitParaIndex = $i
        break
    }
}
$pHabit = $d.Paragraphs($habitParaIndex)
$rngHabit = $d.Range($pHabit.Range.Start, $pHabit.Range.End - 1)
$bodyHabit = '<w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Petra (neutral skeptical): Habit? Who have you been walking home enough for it to be a habit?</w:t></w:r></w:p></w:body>'
$rngHabit.InsertXML((Get-XmlPkg $bodyHabit))

# Re-fetch the (still same-index) paragraph and append a new paragraph after it.
$pHabit = $d.Paragraphs($habitParaIndex)
$pHabit.Range.InsertParagraphAfter()
$pSigh = $d.Paragraphs($habitParaIndex + 1)
# Setting .Text directly first consumes the lone placeholder run that
# InsertParagraphAfter() leaves behind (character-less "orphan" runs can't be
# reliably removed by a zero-length InsertXML call), leaving a single clean run.
$rngSighText = $d.Range($pSigh.Range.Start, $pSigh.Range.End - 1)
$rngSighText.Text = "Petra (neutral sigh):"
# Re-fetch the now-populated range and normalize it via InsertXML so the
# serialized run carries an explicit xml:space="preserve" attribute.
$pSigh2 = $d.Paragraphs($habitParaIndex + 1)
$rngSigh = $d.Range($pSigh2.Range.Start, $pSigh2.Range.End - 1)
$bodySigh = '<w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Petra (neutral sigh):</w:t></w:r></w:p></w:body>'
$rngSigh.InsertXML((Get-XmlPkg $bodySigh))

# ---------------------------------------------------------------------------
# 3) Merge the three runs "Mara: I " + "meant about Lilith" + "." into one
#    run "Mara: I meant about Lilith."
# ---------------------------------------------------------------------------
$maraParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Mara: I meant about Lilith*") {
        $maraParaIndex = $i
        break
    }
}
$pMara = $d.Paragraphs($maraParaIndex)
$rngMara = $d.Range($pMara.Range.Start, $pMara.Range.End - 1)
$bodyMara = '<w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Mara: I meant about Lilith.</w:t></w:r></w:p></w:body>'
$rngMara.InsertXML((Get-XmlPkg $bodyMara))

# ---------------------------------------------------------------------------
# 4) Add xml:space="preserve" (no textual change) on the "Alright." line.
# ---------------------------------------------------------------------------
$alrightParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Alright. I guess we can find a place to sit on our way back*") {
        $alrightParaIndex = $i
        break
    }
}
$pAlright = $d.Paragraphs($alrightParaIndex)
$rngAlright = $d.Range($pAlright.Range.Start, $pAlright.Range.End - 1)
$bodyAlright = '<w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Mara (neutral neutral): Alright. I guess we can find a place to sit on our way back, right?</w:t></w:r></w:p></w:body>'
$rngAlright.InsertXML((Get-XmlPkg $bodyAlright))

# ---------------------------------------------------------------------------
# 5) Add xml:space="preserve" (no textual change) on the "I wouldn't mind
#    talking for a little longer." line.
# ---------------------------------------------------------------------------
$wouldntParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*I wouldn*mind talking for a little longer*") {
        $wouldntParaIndex = $i
        break
    }
}
$pWouldnt = $d.Paragraphs($wouldntParaIndex)
$rngWouldnt = $d.Range($pWouldnt.Range.Start, $pWouldnt.Range.End - 1)
$bodyWouldnt = '<w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Mara (neutral smiling_blushing_eyes): I wouldn’t mind talking for a little longer.</w:t></w:r></w:p></w:body>'
$rngWouldnt.InsertXML((Get-XmlPkg $bodyWouldnt))

Write-Host "All edits applied."
